$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from the end of the paragraph that
#    reads "...φοίτησε στο Σχολείο Δεύτερης Ευκαιρίας Μεσολογγίου." to a new
#    position inside the phrase "...ανταποκρίθηκε [a.antapokrisi] στις
#    απαιτήσεις του προγράμματος σπουδών..." (right after " σ", before
#    "τις").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rngGoBack = $d.Content
$rngGoBack.Find.Execute(" στις απαιτήσεις του προγράμματος σπουδών", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackSplit = $rngGoBack.Start + 2
$goBackRange = $d.Range($goBackSplit, $goBackSplit)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Fill in the placeholder dotted line for the institution name / date:
#      "……………………………………"+"…….."+"………… 20…"+"…."
#    becomes
#      "Ι.Π. Μεσολογγίου"+"…"+" 30/06/2022"+"…"+"…"+"."
# ---------------------------------------------------------------------------
$rngDate = $d.Content
$oldDateText = "……………………………………" + "…….." + "………… 20…" + "…."
$rngDate.Find.Execute($oldDateText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dateBase = $rngDate.Start

$newDateText = "Ι.Π. Μεσολογγίου" + "…" + " 30/06/2022" + "…" + "…" + "."
$rngDate.Text = $newDateText

# Force the replaced text to stay split into separate runs (matching the
# target structure) by nudging formatting on each sub-range back to its
# original value.
$dateSplits = @(16, 17, 28, 29, 30)
$prev = 0
foreach ($off in $dateSplits) {
    $piece = $d.Range($dateBase + $prev, $dateBase + $off)
    $piece.Font.Bold = 1
    $piece.Font.Bold = 0
    $prev = $off
}

# ---------------------------------------------------------------------------
# 3) Fill in the director gender-suffix placeholder:
#      "………… ΔΙΕΥΘΥΝΤ…………"
#    becomes
#      "Η"+" ΔΙΕΥΘΥΝΤ"+"ΡΙΑ"
# ---------------------------------------------------------------------------
$rngDir = $d.Content
$oldDirText = "………… ΔΙΕΥΘΥΝΤ…………"
$rngDir.Find.Execute($oldDirText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dirBase = $rngDir.Start

$newDirText = "Η" + " ΔΙΕΥΘΥΝΤ" + "ΡΙΑ"
$rngDir.Text = $newDirText

$dirSplits = @(1, 10)
$prev = 0
foreach ($off in $dirSplits) {
    $piece = $d.Range($dirBase + $prev, $dirBase + $off)
    $piece.Font.Bold = 1
    $piece.Font.Bold = 0
    $prev = $off
}

Write-Output "ok"
